# Excel COM-interop edit script
# Commit: "Update gh-pages to output generated at 456a3b4"
# The source data refresh bumped the "想去人数" (interest/attendee count)
# figures in column F for a number of rows across all four sheets
# (展览 / 演出 / 本地生活 / 全部类型 -- the last one mirrors rows from the
# first three, merged/sorted by date). Every other cell is untouched.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1367  # was 1362
$ws.Range("F5").Value = 5622  # was 5612
$ws.Range("F6").Value = 412  # was 407
$ws.Range("F7").Value = 1036  # was 1034
$ws.Range("F8").Value = 2743  # was 2721
$ws.Range("F9").Value = 6390  # was 6359
$ws.Range("F10").Value = 174  # was 171
$ws.Range("F11").Value = 1224  # was 1208
$ws.Range("F12").Value = 709  # was 703
$ws.Range("F13").Value = 83  # was 84
$ws.Range("F15").Value = 1099  # was 1096
$ws.Range("F17").Value = 70  # was 64
$ws.Range("F19").Value = 143  # was 142
$ws.Range("F21").Value = 886  # was 875
$ws.Range("F22").Value = 26  # was 24
$ws.Range("F23").Value = 81  # was 79
$ws.Range("F25").Value = 1127  # was 1125
$ws.Range("F28").Value = 220  # was 210
$ws.Range("F29").Value = 42  # was 41
$ws.Range("F30").Value = 222  # was 213
$ws.Range("F31").Value = 1150  # was 1147
$ws.Range("F32").Value = 43  # was 41
$ws.Range("F33").Value = 73  # was 69

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F5").Value = 507  # was 505
$ws.Range("F10").Value = 110  # was 109
$ws.Range("F13").Value = 1  # was 0
$ws.Range("F15").Value = 541  # was 540
$ws.Range("F19").Value = 172  # was 171
$ws.Range("F20").Value = 126  # was 125
$ws.Range("F23").Value = 85  # was 83
$ws.Range("F24").Value = 43  # was 44
$ws.Range("F26").Value = 107  # was 106
$ws.Range("F27").Value = 600  # was 573
$ws.Range("F28").Value = 937  # was 935
$ws.Range("F29").Value = 549  # was 545
$ws.Range("F33").Value = 6  # was 5
$ws.Range("F34").Value = 91  # was 89
$ws.Range("F35").Value = 120  # was 119
$ws.Range("F37").Value = 46  # was 45

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F5").Value = 815  # was 816
$ws.Range("F6").Value = 514  # was 513
$ws.Range("F7").Value = 276  # was 274

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F4").Value = 1367  # was 1362
$ws.Range("F6").Value = 815  # was 816
$ws.Range("F9").Value = 514  # was 513
$ws.Range("F10").Value = 276  # was 274
$ws.Range("F11").Value = 276  # was 274
$ws.Range("F12").Value = 507  # was 505
$ws.Range("F14").Value = 5622  # was 5612
$ws.Range("F15").Value = 412  # was 407
$ws.Range("F16").Value = 1036  # was 1034
$ws.Range("F17").Value = 2743  # was 2721
$ws.Range("F19").Value = 6390  # was 6359
$ws.Range("F20").Value = 110  # was 109
$ws.Range("F21").Value = 174  # was 171
$ws.Range("F22").Value = 1224  # was 1208
$ws.Range("F24").Value = 541  # was 540
$ws.Range("F25").Value = 709  # was 703
$ws.Range("F26").Value = 83  # was 84
$ws.Range("F27").Value = 1099  # was 1096
$ws.Range("F28").Value = 126  # was 125
$ws.Range("F29").Value = 70  # was 64
$ws.Range("F30").Value = 143  # was 142
$ws.Range("F32").Value = 886  # was 875
$ws.Range("F33").Value = 85  # was 83
$ws.Range("F34").Value = 81  # was 79
$ws.Range("F35").Value = 1127  # was 1125
$ws.Range("F37").Value = 107  # was 106
$ws.Range("F39").Value = 937  # was 935
$ws.Range("F40").Value = 549  # was 545
$ws.Range("F41").Value = 220  # was 210
$ws.Range("F42").Value = 42  # was 41
$ws.Range("F44").Value = 222  # was 213
$ws.Range("F45").Value = 6  # was 5
$ws.Range("F46").Value = 91  # was 89
$ws.Range("F47").Value = 120  # was 119
$ws.Range("F49").Value = 73  # was 69
$ws.Range("F50").Value = 46  # was 45
